$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(213).Insert()

$ws.Range("A213").Value = 10
$ws.Range("B213").Value = 'Vega Modelo de Temuco'
$ws.Range("C213").Value = 'La Araucanía'
$ws.Range("D213").Value = 44582
$ws.Range("E213").Value = 9
$ws.Range("F213").Value = 100114014
$ws.Range("G213").Value = 'Betarraga'
$ws.Range("H213").Value = 'Sin especificar'
$ws.Range("I213").Value = 'Primera'
$ws.Range("J213").Value = 70
$ws.Range("K213").Value = 7000
$ws.Range("L213").Value = 8000
$ws.Range("M213").Value = 7429
$ws.Range("N213").Value = '$/docena de paquetes'
$ws.Range("O213").Value = 'Provincia de Cautín'
$ws.Range("P213").Value = 619
$ws.Range("Q213").Value = 12
$ws.Range("R213").Value = 'Hortaliza'
